# cryptos.xlsx symbol-list refresh — GitHub Actions scraper run
# Updates the "Price" (D) and "Volume(1h)" (E) columns for each coin row
# with newly scraped quotes. Only the cells explicitly listed below change;
# everything else (labels, links, dates, styles) is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

function Set-TextCell([string]$addr, [string]$val) {
    # The sheet stores these numbers/percentages as plain text (inline
    # strings), not numeric cells. Writing a string like "277.16" or
    # "1.66%" straight into .Value would get auto-coerced by Excel into a
    # real number/percentage (and swap in a numeric style). Force the
    # cell to Text format first, assign the literal text, then restore
    # the original style so formatting is unaffected.
    $cell = $ws.Range($addr)
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $val
    $cell.Style = $origStyle
}

Set-TextCell "D2" "277.16"
Set-TextCell "E2" "1.66%"
Set-TextCell "D3" "27.15"
Set-TextCell "E3" "1.15%"
Set-TextCell "D4" "4.937"
Set-TextCell "E4" "0.80%"
Set-TextCell "D5" "0.06410"
Set-TextCell "E5" "1.39%"
Set-TextCell "D6" "6.922"
Set-TextCell "E6" "0.16%"
Set-TextCell "E7" "-6.62%"
Set-TextCell "D8" "0.8820"
Set-TextCell "E8" "-0.13%"
Set-TextCell "D9" "0.1525"
Set-TextCell "E9" "4.35%"
Set-TextCell "D10" "0.05071"
Set-TextCell "E10" "-0.62%"
Set-TextCell "D11" "0.07523"
Set-TextCell "E11" "1.56%"
Set-TextCell "D12" "0.02878"
Set-TextCell "E12" "-8.61%"
Set-TextCell "D13" "0.09014"
Set-TextCell "E13" "-0.20%"
Set-TextCell "D14" "0.001568"
Set-TextCell "E14" "0.02%"
Set-TextCell "D15" "0.0006442"
Set-TextCell "E15" "2.20%"
Set-TextCell "D16" "0.005982"
Set-TextCell "E16" "-0.43%"
Set-TextCell "E17" "-0.30%"
Set-TextCell "E18" "-0.89%"
Set-TextCell "D19" "2.272"
Set-TextCell "E19" "-0.01%"
Set-TextCell "E20" "0.21%"
Set-TextCell "E21" "0.36%"
Set-TextCell "D22" "3.908"
Set-TextCell "E22" "-0.08%"
Set-TextCell "D23" "0.04435"
Set-TextCell "E23" "2.13%"
Set-TextCell "D24" "0.001175"
Set-TextCell "E24" "-0.07%"
Set-TextCell "D25" "0.003875"
Set-TextCell "E25" "6.10%"
Set-TextCell "D26" "0.0001202"
Set-TextCell "E27" "14.12%"
Set-TextCell "D40" "0.04142"
Set-TextCell "E40" "2.41%"
Set-TextCell "D41" "0.006793"
Set-TextCell "E41" "2.56%"
Set-TextCell "E42" "1.25%"
Set-TextCell "D43" "0.002293"
Set-TextCell "E43" "7.66%"
Set-TextCell "D44" "0.01124"
Set-TextCell "E44" "-10.52%"
Set-TextCell "E45" "-2.29%"
Set-TextCell "D46" "1.482"
Set-TextCell "E46" "-37.06%"
Set-TextCell "D47" "0.02025"
Set-TextCell "E47" "-4.41%"
